$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the full target range as Text so that numeric-looking
# values (totalRuns, totalBalls, total4s, total6s, sr) are stored as
# text strings (t="str"/shared-string), matching the source data,
# instead of being auto-converted to numbers by Excel.
$ws.Range("A1:K16").NumberFormat = "@"

# --- Header row (row 1) ---
$ws.Cells.Item(1, 1).Value = 'venue'
$ws.Cells.Item(1, 2).Value = 'date'
$ws.Cells.Item(1, 3).Value = 'result'
$ws.Cells.Item(1, 4).Value = 'ownTeam'
$ws.Cells.Item(1, 5).Value = 'oppTeam'
$ws.Cells.Item(1, 6).Value = 'batsman'
$ws.Cells.Item(1, 7).Value = 'totalRuns'
$ws.Cells.Item(1, 8).Value = 'totalBalls'
$ws.Cells.Item(1, 9).Value = 'total4s'
$ws.Cells.Item(1, 10).Value = 'total6s'
$ws.Cells.Item(1, 11).Value = 'sr'

# --- Data rows (rows 2-16) ---
# row 2
$ws.Cells.Item(2, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(2, 2).Value = ' October 25 2020'
$ws.Cells.Item(2, 3).Value = 'Royals won by 8 wickets (with 10 balls remaining)'
$ws.Cells.Item(2, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(2, 5).Value = 'Rajasthan Royals'
$ws.Cells.Item(2, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(2, 7).Value = '40'
$ws.Cells.Item(2, 8).Value = '26'
$ws.Cells.Item(2, 9).Value = '4'
$ws.Cells.Item(2, 10).Value = '1'
$ws.Cells.Item(2, 11).Value = '153.84'

# row 3
$ws.Cells.Item(3, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(3, 2).Value = ' October 28 2020'
$ws.Cells.Item(3, 3).Value = 'Mumbai won by 5 wickets (with 5 balls remaining)'
$ws.Cells.Item(3, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(3, 5).Value = 'Royal Challengers Bangalore'
$ws.Cells.Item(3, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(3, 7).Value = '79'
$ws.Cells.Item(3, 8).Value = '43'
$ws.Cells.Item(3, 9).Value = '10'
$ws.Cells.Item(3, 10).Value = '3'
$ws.Cells.Item(3, 11).Value = '183.72'

# row 4
$ws.Cells.Item(4, 1).Value = ' Sharjah'
$ws.Cells.Item(4, 2).Value = ' November 03 2020'
$ws.Cells.Item(4, 3).Value = 'Sunrisers won by 10 wickets (with 17 balls remaining)'
$ws.Cells.Item(4, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(4, 5).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(4, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(4, 7).Value = '36'
$ws.Cells.Item(4, 8).Value = '29'
$ws.Cells.Item(4, 9).Value = '5'
$ws.Cells.Item(4, 10).Value = '0'
$ws.Cells.Item(4, 11).Value = '124.13'

# row 5
$ws.Cells.Item(5, 1).Value = ' Dubai (DSC)'
$ws.Cells.Item(5, 2).Value = ' November 05 2020'
$ws.Cells.Item(5, 3).Value = 'Mumbai won by 57 runs'
$ws.Cells.Item(5, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(5, 5).Value = 'Delhi Capitals'
$ws.Cells.Item(5, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(5, 7).Value = '51'
$ws.Cells.Item(5, 8).Value = '38'
$ws.Cells.Item(5, 9).Value = '6'
$ws.Cells.Item(5, 10).Value = '2'
$ws.Cells.Item(5, 11).Value = '134.21'

# row 6
$ws.Cells.Item(6, 1).Value = ' Dubai (DSC)'
$ws.Cells.Item(6, 2).Value = ' November 10 2020'
$ws.Cells.Item(6, 3).Value = 'Mumbai won by 5 wickets (with 8 balls remaining)'
$ws.Cells.Item(6, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(6, 5).Value = 'Delhi Capitals'
$ws.Cells.Item(6, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(6, 7).Value = '19'
$ws.Cells.Item(6, 8).Value = '20'
$ws.Cells.Item(6, 9).Value = '1'
$ws.Cells.Item(6, 10).Value = '1'
$ws.Cells.Item(6, 11).Value = '95.00'

# row 7
$ws.Cells.Item(7, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(7, 2).Value = ' October 16 2020'
$ws.Cells.Item(7, 3).Value = 'Mumbai won by 8 wickets (with 19 balls remaining)'
$ws.Cells.Item(7, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(7, 5).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(7, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(7, 7).Value = '10'
$ws.Cells.Item(7, 8).Value = '10'
$ws.Cells.Item(7, 9).Value = '1'
$ws.Cells.Item(7, 10).Value = '0'
$ws.Cells.Item(7, 11).Value = '100.00'

# row 8
$ws.Cells.Item(8, 1).Value = ' Dubai (DSC)'
$ws.Cells.Item(8, 2).Value = ' October 31 2020'
$ws.Cells.Item(8, 3).Value = 'Mumbai won by 9 wickets (with 34 balls remaining)'
$ws.Cells.Item(8, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(8, 5).Value = 'Delhi Capitals'
$ws.Cells.Item(8, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(8, 7).Value = '12'
$ws.Cells.Item(8, 8).Value = '11'
$ws.Cells.Item(8, 9).Value = '1'
$ws.Cells.Item(8, 10).Value = '0'
$ws.Cells.Item(8, 11).Value = '109.09'

# row 9
$ws.Cells.Item(9, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(9, 2).Value = ' October 01 2020'
$ws.Cells.Item(9, 3).Value = 'Mumbai won by 48 runs'
$ws.Cells.Item(9, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(9, 5).Value = 'Kings XI Punjab'
$ws.Cells.Item(9, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(9, 7).Value = '10'
$ws.Cells.Item(9, 8).Value = '7'
$ws.Cells.Item(9, 9).Value = '2'
$ws.Cells.Item(9, 10).Value = '0'
$ws.Cells.Item(9, 11).Value = '142.85'

# row 10
$ws.Cells.Item(10, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(10, 2).Value = ' September 19 2020'
$ws.Cells.Item(10, 3).Value = 'Super Kings won by 5 wickets (with 4 balls remaining)'
$ws.Cells.Item(10, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(10, 5).Value = 'Chennai Super Kings'
$ws.Cells.Item(10, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(10, 7).Value = '17'
$ws.Cells.Item(10, 8).Value = '16'
$ws.Cells.Item(10, 9).Value = '2'
$ws.Cells.Item(10, 10).Value = '0'
$ws.Cells.Item(10, 11).Value = '106.25'

# row 11
$ws.Cells.Item(11, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(11, 2).Value = ' October 06 2020'
$ws.Cells.Item(11, 3).Value = 'Mumbai won by 57 runs'
$ws.Cells.Item(11, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(11, 5).Value = 'Rajasthan Royals'
$ws.Cells.Item(11, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(11, 7).Value = '79'
$ws.Cells.Item(11, 8).Value = '47'
$ws.Cells.Item(11, 9).Value = '11'
$ws.Cells.Item(11, 10).Value = '2'
$ws.Cells.Item(11, 11).Value = '168.08'

# row 12
$ws.Cells.Item(12, 1).Value = ' Dubai (DSC)'
$ws.Cells.Item(12, 2).Value = ' October 18 2020'
$ws.Cells.Item(12, 3).Value = 'Match tied (Kings XI won the one-over eliminator)'
$ws.Cells.Item(12, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(12, 5).Value = 'Kings XI Punjab'
$ws.Cells.Item(12, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(12, 7).Value = '0'
$ws.Cells.Item(12, 8).Value = '4'
$ws.Cells.Item(12, 9).Value = '0'
$ws.Cells.Item(12, 10).Value = '0'
$ws.Cells.Item(12, 11).Value = '0.00'

# row 13
$ws.Cells.Item(13, 1).Value = ' Dubai (DSC)'
$ws.Cells.Item(13, 2).Value = ' September 28 2020'
$ws.Cells.Item(13, 3).Value = 'Match tied (RCB won the one-over eliminator)'
$ws.Cells.Item(13, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(13, 5).Value = 'Royal Challengers Bangalore'
$ws.Cells.Item(13, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(13, 7).Value = '0'
$ws.Cells.Item(13, 8).Value = '2'
$ws.Cells.Item(13, 9).Value = '0'
$ws.Cells.Item(13, 10).Value = '0'
$ws.Cells.Item(13, 11).Value = '0.00'

# row 14
$ws.Cells.Item(14, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(14, 2).Value = ' September 23 2020'
$ws.Cells.Item(14, 3).Value = 'Mumbai won by 49 runs'
$ws.Cells.Item(14, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(14, 5).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(14, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(14, 7).Value = '47'
$ws.Cells.Item(14, 8).Value = '28'
$ws.Cells.Item(14, 9).Value = '6'
$ws.Cells.Item(14, 10).Value = '1'
$ws.Cells.Item(14, 11).Value = '167.85'

# row 15
$ws.Cells.Item(15, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(15, 2).Value = ' October 11 2020'
$ws.Cells.Item(15, 3).Value = 'Mumbai won by 5 wickets (with 2 balls remaining)'
$ws.Cells.Item(15, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(15, 5).Value = 'Delhi Capitals'
$ws.Cells.Item(15, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(15, 7).Value = '53'
$ws.Cells.Item(15, 8).Value = '32'
$ws.Cells.Item(15, 9).Value = '6'
$ws.Cells.Item(15, 10).Value = '1'
$ws.Cells.Item(15, 11).Value = '165.62'

# row 16
$ws.Cells.Item(16, 1).Value = ' Sharjah'
$ws.Cells.Item(16, 2).Value = ' October 04 2020'
$ws.Cells.Item(16, 3).Value = 'Mumbai won by 34 runs'
$ws.Cells.Item(16, 4).Value = 'Mumbai Indians'
$ws.Cells.Item(16, 5).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(16, 6).Value = 'Suryakumar Yadav '
$ws.Cells.Item(16, 7).Value = '27'
$ws.Cells.Item(16, 8).Value = '18'
$ws.Cells.Item(16, 9).Value = '6'
$ws.Cells.Item(16, 10).Value = '0'
$ws.Cells.Item(16, 11).Value = '150.00'

Write-Output "Updated sheet: wrote 11 headers + 15 data rows (A1:K16)."
